$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add the new "Connector" column (F) to Sheet1 ---
$ws1.Range("F1").Value = "Connector"
$ws1.Range("F2").Value = "AND"
$ws1.Range("F3").Value = "OR"
$ws1.Columns.Item(6).AutoFit()

# --- Duplicate Sheet1 as the template for the new "Country_Group_A" scenario ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# --- Overwrite the scenario-specific cells on the new sheet ---
$ws2.Range("A2").Value = "Country_Group_A"
$ws2.Range("B2").Value = "Countries name starting with A"
$ws2.Range("C2").Value = "Phone"
$ws2.Range("D2").Value = "contains"
$ws2.Range("E2").Value = "z"

$ws2.Range("C3").Value = "Report Group"
$ws2.Range("D3").Value = "does not end with"
$ws2.Range("E3").Value = "y"

$ws2.Range("C4").Value = "Country"
$ws2.Range("D4").Value = "starts with"
$ws2.Range("E4").Value = "a"

# Re-fit the columns whose content actually changed
$ws2.Columns.Item(3).AutoFit()
$ws2.Columns.Item(4).AutoFit()
$ws2.Columns.Item(6).AutoFit()

# --- Selections: both sheets end up with F1:F3 selected (Connector verification) ---
$ws1.Range("F1:F3").Select()

$ws2.Activate()
$ws2.Range("F1:F3").Select()
